$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Range("E4").Value = 21
$ws.Range("F4").Value = 10
$ws.Range("H4").Value = 4
$ws.Range("K4").Value = 61
$ws.Range("L4").Value = 42
$ws.Range("M4").Value = 35
$ws.Range("N4").Value = 37
$ws.Range("Q4").Value = 5
$ws.Range("R4").Value = 4.8
$ws.Range("S4").Value = 252
$ws.Range("T4").Value = 12
$ws.Range("AB4").Value = 298
$ws.Range("AC4").Value = 1194
$ws.Range("AD4").Value = 59
$ws.Range("AJ4").Value = 1.492
$ws.Range("AK4").Value = 0.288
$ws.Range("AM4").Value = 1.84

# Row 10
$ws.Range("E10").Value = 33
$ws.Range("J10").Value = 49
$ws.Range("K10").Value = 47
$ws.Range("L10").Value = 28
$ws.Range("M10").Value = 25
$ws.Range("N10").Value = 16
$ws.Range("O10").Value = 40
$ws.Range("R10").Value = 4.59
$ws.Range("S10").Value = 187
$ws.Range("T10").Value = 7
$ws.Range("W10").Value = 3
$ws.Range("AB10").Value = 208
$ws.Range("AC10").Value = 784
$ws.Range("AH10").Value = 16
$ws.Range("AJ10").Value = 1.286
$ws.Range("AK10").Value = 0.294
$ws.Range("AL10").Value = 7.3
$ws.Range("AM10").Value = 2.5

# Row 11
$ws.Range("E11").Value = 13
$ws.Range("F11").Value = 13
$ws.Range("G11").Value = 4
$ws.Range("J11").Value = 76
$ws.Range("K11").Value = 55
$ws.Range("L11").Value = 22
$ws.Range("M11").Value = 22
$ws.Range("O11").Value = 65
$ws.Range("P11").Value = 9
$ws.Range("R11").Value = 2.61
$ws.Range("S11").Value = 278
$ws.Range("T11").Value = 8
$ws.Range("AB11").Value = 300
$ws.Range("AC11").Value = 1157
$ws.Range("AE11").Value = 17
$ws.Range("AJ11").Value = 1.013
$ws.Range("AK11").Value = 0.225
$ws.Range("AL11").Value = 7.7
$ws.Range("AM11").Value = 2.95

# Row 13
$ws.Range("E13").Value = 11
$ws.Range("F13").Value = 10
$ws.Range("H13").Value = 4
$ws.Range("J13").Value = 53.1
$ws.Range("K13").Value = 49
$ws.Range("L13").Value = 35
$ws.Range("M13").Value = 32
$ws.Range("N13").Value = 19
$ws.Range("O13").Value = 53
$ws.Range("P13").Value = 10
$ws.Range("R13").Value = 5.4
$ws.Range("S13").Value = 204
$ws.Range("T13").Value = 7
$ws.Range("AB13").Value = 229
$ws.Range("AC13").Value = 882
$ws.Range("AD13").Value = 64
$ws.Range("AG13").Value = 37
$ws.Range("AJ13").Value = 1.275
$ws.Range("AK13").Value = 0.275
$ws.Range("AL13").Value = 8.9
$ws.Range("AM13").Value = 2.79

# Row 19
$ws.Range("E19").Value = 5
$ws.Range("F19").Value = 4
$ws.Range("G19").Value = 2
$ws.Range("J19").Value = 28.1
$ws.Range("K19").Value = 26
$ws.Range("L19").Value = 12
$ws.Range("M19").Value = 12
$ws.Range("O19").Value = 17
$ws.Range("R19").Value = 3.81
$ws.Range("S19").Value = 107
$ws.Range("W19").Value = 2
$ws.Range("Z19").Value = 1
$ws.Range("AB19").Value = 112
$ws.Range("AC19").Value = 425
$ws.Range("AD19").Value = 62
$ws.Range("AE19").Value = 15
$ws.Range("AF19").Value = 8
$ws.Range("AG19").Value = 38
$ws.Range("AH19").Value = 32
$ws.Range("AI19").Value = 5
$ws.Range("AJ19").Value = 1.024
$ws.Range("AK19").Value = 0.273
$ws.Range("AL19").Value = 5.4
$ws.Range("AM19").Value = 5.67

# Row 25
$ws.Range("E25").Value = 21
$ws.Range("F25").Value = 21
$ws.Range("G25").Value = 5
$ws.Range("J25").Value = 127
$ws.Range("K25").Value = 85
$ws.Range("N25").Value = 31
$ws.Range("O25").Value = 137
$ws.Range("Q25").Value = 5
$ws.Range("R25").Value = 1.91
$ws.Range("S25").Value = 456
$ws.Range("T25").Value = 23
$ws.Range("AB25").Value = 492
$ws.Range("AC25").Value = 1987
$ws.Range("AI25").Value = 7
$ws.Range("AJ25").Value = 0.913
$ws.Range("AK25").Value = 0.252
$ws.Range("AM25").Value = 4.42

# Row 32
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 16
$ws.Range("G32").Value = 4
$ws.Range("J32").Value = 84.2
$ws.Range("K32").Value = 100
$ws.Range("L32").Value = 54
$ws.Range("M32").Value = 53
$ws.Range("N32").Value = 31
$ws.Range("O32").Value = 70
$ws.Range("R32").Value = 5.63
$ws.Range("S32").Value = 342
$ws.Range("T32").Value = 24
$ws.Range("W32").Value = 6
$ws.Range("X32").Value = 3
$ws.Range("AB32").Value = 381
$ws.Range("AC32").Value = 1434
$ws.Range("AE32").Value = 17
$ws.Range("AG32").Value = 50
$ws.Range("AJ32").Value = 1.547
$ws.Range("AK32").Value = 0.324
$ws.Range("AL32").Value = 7.4
$ws.Range("AM32").Value = 2.26

# Row 63
$ws.Range("E63").Value = 17
$ws.Range("F63").Value = 17
$ws.Range("G63").Value = 6
$ws.Range("J63").Value = 86.09999999999999
$ws.Range("K63").Value = 75
$ws.Range("L63").Value = 41
$ws.Range("M63").Value = 41
$ws.Range("N63").Value = 42
$ws.Range("O63").Value = 76
$ws.Range("P63").Value = 11
$ws.Range("Q63").Value = 3
$ws.Range("R63").Value = 4.27
$ws.Range("S63").Value = 327
$ws.Range("T63").Value = 15
$ws.Range("AB63").Value = 373
$ws.Range("AC63").Value = 1478
$ws.Range("AD63").Value = 62
$ws.Range("AE63").Value = 15
$ws.Range("AG63").Value = 38
$ws.Range("AH63").Value = 23
$ws.Range("AI63").Value = 7
$ws.Range("AJ63").Value = 1.355
$ws.Range("AK63").Value = 0.266
$ws.Range("AL63").Value = 7.9
$ws.Range("AM63").Value = 1.81

# Row 72
$ws.Range("E72").Value = 20
$ws.Range("F72").Value = 20
$ws.Range("J72").Value = 99
$ws.Range("K72").Value = 93
$ws.Range("L72").Value = 42
$ws.Range("M72").Value = 41
$ws.Range("N72").Value = 47
$ws.Range("P72").Value = 13
$ws.Range("R72").Value = 3.73
$ws.Range("S72").Value = 374
$ws.Range("W72").Value = 8
$ws.Range("AB72").Value = 431
$ws.Range("AC72").Value = 1620
$ws.Range("AJ72").Value = 1.414
$ws.Range("AK72").Value = 0.254
$ws.Range("AL72").Value = 4.5
$ws.Range("AM72").Value = 1.06

# Row 76
$ws.Range("E76").Value = 20
$ws.Range("F76").Value = 17
$ws.Range("J76").Value = 101
$ws.Range("K76").Value = 96
$ws.Range("L76").Value = 57
$ws.Range("N76").Value = 47
$ws.Range("O76").Value = 85
$ws.Range("R76").Value = 4.19
$ws.Range("S76").Value = 387
$ws.Range("T76").Value = 22
$ws.Range("Y76").Value = 13
$ws.Range("AB76").Value = 444
$ws.Range("AC76").Value = 1728
$ws.Range("AD76").Value = 61
$ws.Range("AG76").Value = 35
$ws.Range("AI76").Value = 8
$ws.Range("AJ76").Value = 1.416
$ws.Range("AL76").Value = 7.6
$ws.Range("AM76").Value = 1.81

# Row 78
$ws.Range("E78").Value = 28
$ws.Range("J78").Value = 32.1
$ws.Range("K78").Value = 29
$ws.Range("L78").Value = 20
$ws.Range("M78").Value = 14
$ws.Range("N78").Value = 17
$ws.Range("O78").Value = 26
$ws.Range("R78").Value = 3.9
$ws.Range("S78").Value = 123
$ws.Range("T78").Value = 8
$ws.Range("X78").Value = 3
$ws.Range("AB78").Value = 145
$ws.Range("AC78").Value = 565
$ws.Range("AG78").Value = 57
$ws.Range("AH78").Value = 20
$ws.Range("AJ78").Value = 1.423
$ws.Range("AK78").Value = 0.276
$ws.Range("AL78").Value = 7.2
$ws.Range("AM78").Value = 1.53

# Row 85
$ws.Range("E85").Value = 13
$ws.Range("J85").Value = 24
$ws.Range("N85").Value = 11
$ws.Range("O85").Value = 20
$ws.Range("R85").Value = 2.63
$ws.Range("S85").Value = 93
$ws.Range("AB85").Value = 104
$ws.Range("AC85").Value = 401
$ws.Range("AF85").Value = 11
$ws.Range("AG85").Value = 53
$ws.Range("AJ85").Value = 1.458
$ws.Range("AK85").Value = 0.319
$ws.Range("AL85").Value = 7.5
$ws.Range("AM85").Value = 1.82

# Row 89
$ws.Range("E89").Value = 17
$ws.Range("H89").Value = 7
$ws.Range("J89").Value = 83.09999999999999
$ws.Range("K89").Value = 101
$ws.Range("L89").Value = 62
$ws.Range("M89").Value = 60
$ws.Range("O89").Value = 95
$ws.Range("P89").Value = 15
$ws.Range("R89").Value = 6.48
$ws.Range("S89").Value = 346
$ws.Range("T89").Value = 13
$ws.Range("Y89").Value = 8
$ws.Range("AB89").Value = 378
$ws.Range("AC89").Value = 1399
$ws.Range("AF89").Value = 12
$ws.Range("AJ89").Value = 1.536
$ws.Range("AK89").Value = 0.363
$ws.Range("AM89").Value = 3.52

# Row 92
$ws.Range("E92").Value = 37
$ws.Range("J92").Value = 47.1
$ws.Range("O92").Value = 61
$ws.Range("R92").Value = 3.23
$ws.Range("S92").Value = 171
$ws.Range("AB92").Value = 196
$ws.Range("AC92").Value = 828
$ws.Range("AF92").Value = 15
$ws.Range("AJ92").Value = 1.141
$ws.Range("AL92").Value = 11.6
$ws.Range("AM92").Value = 3.21

# Row 93
$ws.Range("E93").Value = 12
$ws.Range("F93").Value = 11
$ws.Range("G93").Value = 4
$ws.Range("J93").Value = 57.2
$ws.Range("K93").Value = 72
$ws.Range("L93").Value = 31
$ws.Range("M93").Value = 29
$ws.Range("N93").Value = 11
$ws.Range("O93").Value = 48
$ws.Range("P93").Value = 11
$ws.Range("R93").Value = 4.53
$ws.Range("S93").Value = 238
$ws.Range("T93").Value = 15
$ws.Range("X93").Value = 2
$ws.Range("Y93").Value = 2
$ws.Range("AB93").Value = 254
$ws.Range("AC93").Value = 950
$ws.Range("AG93").Value = 35
$ws.Range("AJ93").Value = 1.439
$ws.Range("AK93").Value = 0.337
$ws.Range("AL93").Value = 7.5
$ws.Range("AM93").Value = 4.36

# Row 94
$ws.Range("E94").Value = 20
$ws.Range("F94").Value = 20
$ws.Range("G94").Value = 6
$ws.Range("J94").Value = 116
$ws.Range("K94").Value = 108
$ws.Range("L94").Value = 60
$ws.Range("M94").Value = 55
$ws.Range("O94").Value = 99
$ws.Range("P94").Value = 18
$ws.Range("R94").Value = 4.27
$ws.Range("S94").Value = 439
$ws.Range("W94").Value = 8
$ws.Range("Y94").Value = 15
$ws.Range("AB94").Value = 483
$ws.Range("AC94").Value = 1904
$ws.Range("AD94").Value = 64
$ws.Range("AH94").Value = 19
$ws.Range("AI94").Value = 8
$ws.Range("AJ94").Value = 1.216
$ws.Range("AK94").Value = 0.278
$ws.Range("AL94").Value = 7.7
$ws.Range("AM94").Value = 3

# Row 103
$ws.Range("E103").Value = 15
$ws.Range("F103").Value = 15
$ws.Range("J103").Value = 80.2
$ws.Range("K103").Value = 59
$ws.Range("L103").Value = 19
$ws.Range("M103").Value = 16
$ws.Range("N103").Value = 36
$ws.Range("O103").Value = 79
$ws.Range("P103").Value = 5
$ws.Range("R103").Value = 1.79
$ws.Range("S103").Value = 290
$ws.Range("T103").Value = 14
$ws.Range("AB103").Value = 327
$ws.Range("AC103").Value = 1316
$ws.Range("AD103").Value = 60
$ws.Range("AJ103").Value = 1.178
$ws.Range("AK103").Value = 0.262
$ws.Range("AL103").Value = 8.800000000000001
$ws.Range("AM103").Value = 2.19

# Row 105
$ws.Range("E105").Value = 42
$ws.Range("J105").Value = 45.1
$ws.Range("K105").Value = 44
$ws.Range("L105").Value = 22
$ws.Range("M105").Value = 21
$ws.Range("O105").Value = 55
$ws.Range("Q105").Value = 2
$ws.Range("R105").Value = 4.17
$ws.Range("S105").Value = 172
$ws.Range("Y105").Value = 3
$ws.Range("AB105").Value = 193
$ws.Range("AC105").Value = 792
$ws.Range("AG105").Value = 49
$ws.Range("AH105").Value = 26
$ws.Range("AI105").Value = 3
$ws.Range("AJ105").Value = 1.39
$ws.Range("AK105").Value = 0.354
$ws.Range("AL105").Value = 10.9
$ws.Range("AM105").Value = 2.89

# Row 113
$ws.Range("E113").Value = 11
$ws.Range("F113").Value = 11
$ws.Range("H113").Value = 5
$ws.Range("J113").Value = 60
$ws.Range("K113").Value = 55
$ws.Range("L113").Value = 31
$ws.Range("M113").Value = 31
$ws.Range("N113").Value = 13
$ws.Range("O113").Value = 59
$ws.Range("R113").Value = 4.65
$ws.Range("S113").Value = 226
$ws.Range("T113").Value = 8
$ws.Range("X113").Value = 4
$ws.Range("Z113").Value = 2
$ws.Range("AB113").Value = 246
$ws.Range("AC113").Value = 977
$ws.Range("AG113").Value = 46
$ws.Range("AH113").Value = 25
$ws.Range("AI113").Value = 4
$ws.Range("AJ113").Value = 1.133
$ws.Range("AK113").Value = 0.288
$ws.Range("AL113").Value = 8.9
$ws.Range("AM113").Value = 4.54

# Row 122
$ws.Range("E122").Value = 20
$ws.Range("F122").Value = 20
$ws.Range("H122").Value = 6
$ws.Range("J122").Value = 114
$ws.Range("K122").Value = 106
$ws.Range("L122").Value = 59
$ws.Range("M122").Value = 59
$ws.Range("N122").Value = 42
$ws.Range("O122").Value = 111
$ws.Range("P122").Value = 19
$ws.Range("R122").Value = 4.66
$ws.Range("S122").Value = 430
$ws.Range("T122").Value = 25
$ws.Range("AB122").Value = 480
$ws.Range("AC122").Value = 1834
$ws.Range("AG122").Value = 48
$ws.Range("AH122").Value = 21
$ws.Range("AJ122").Value = 1.298
$ws.Range("AK122").Value = 0.289
$ws.Range("AM122").Value = 2.64

# Row 131
$ws.Range("E131").Value = 37
$ws.Range("J131").Value = 42
$ws.Range("K131").Value = 42
$ws.Range("L131").Value = 22
$ws.Range("M131").Value = 17
$ws.Range("N131").Value = 18
$ws.Range("R131").Value = 3.64
$ws.Range("S131").Value = 161
$ws.Range("T131").Value = 13
$ws.Range("AB131").Value = 182
$ws.Range("AC131").Value = 663
$ws.Range("AG131").Value = 59
$ws.Range("AH131").Value = 24
$ws.Range("AJ131").Value = 1.429
$ws.Range("AK131").Value = 0.342
$ws.Range("AL131").Value = 9
$ws.Range("AM131").Value = 2.33

# Row 134
$ws.Range("E134").Value = 21
$ws.Range("F134").Value = 21
$ws.Range("H134").Value = 11
$ws.Range("J134").Value = 121
$ws.Range("K134").Value = 120
$ws.Range("L134").Value = 80
$ws.Range("M134").Value = 75
$ws.Range("N134").Value = 46
$ws.Range("O134").Value = 116
$ws.Range("P134").Value = 23
$ws.Range("R134").Value = 5.58
$ws.Range("S134").Value = 470
$ws.Range("T134").Value = 23
$ws.Range("AB134").Value = 522
$ws.Range("AC134").Value = 1981
$ws.Range("AH134").Value = 22
$ws.Range("AJ134").Value = 1.372
$ws.Range("AM134").Value = 2.52

# Row 143
$ws.Range("E143").Value = 10
$ws.Range("J143").Value = 19
$ws.Range("O143").Value = 11
$ws.Range("R143").Value = 5.68
$ws.Range("S143").Value = 70
$ws.Range("AB143").Value = 79
$ws.Range("AC143").Value = 288
$ws.Range("AD143").Value = 66
$ws.Range("AF143").Value = 10
$ws.Range("AG143").Value = 48
$ws.Range("AJ143").Value = 1.263
$ws.Range("AK143").Value = 0.254
$ws.Range("AL143").Value = 5.2
$ws.Range("AM143").Value = 1.57

# Row 144
$ws.Range("E144").Value = 20
$ws.Range("F144").Value = 20
$ws.Range("H144").Value = 8
$ws.Range("J144").Value = 102.1
$ws.Range("K144").Value = 100
$ws.Range("L144").Value = 58
$ws.Range("M144").Value = 55
$ws.Range("N144").Value = 42
$ws.Range("O144").Value = 93
$ws.Range("P144").Value = 12
$ws.Range("Q144").Value = 6
$ws.Range("R144").Value = 4.84
$ws.Range("S144").Value = 390
$ws.Range("T144").Value = 25
$ws.Range("U144").Value = 4
$ws.Range("AB144").Value = 444
$ws.Range("AC144").Value = 1767
$ws.Range("AJ144").Value = 1.388
$ws.Range("AK144").Value = 0.303
$ws.Range("AL144").Value = 8.199999999999999
$ws.Range("AM144").Value = 2.21

# Row 145
$ws.Range("E145").Value = 21
$ws.Range("F145").Value = 21
$ws.Range("J145").Value = 119.2
$ws.Range("K145").Value = 118
$ws.Range("L145").Value = 68
$ws.Range("M145").Value = 64
$ws.Range("N145").Value = 41
$ws.Range("O145").Value = 84
$ws.Range("R145").Value = 4.81
$ws.Range("S145").Value = 468
$ws.Range("T145").Value = 24
$ws.Range("W145").Value = 10
$ws.Range("AB145").Value = 513
$ws.Range("AC145").Value = 1933
$ws.Range("AJ145").Value = 1.329
$ws.Range("AK145").Value = 0.265
$ws.Range("AL145").Value = 6.3
$ws.Range("AM145").Value = 2.05

# Row 151
$ws.Range("E151").Value = 19
$ws.Range("J151").Value = 39.1
$ws.Range("O151").Value = 20
$ws.Range("R151").Value = 2.97
$ws.Range("S151").Value = 154
$ws.Range("AB151").Value = 169
$ws.Range("AC151").Value = 592
$ws.Range("AG151").Value = 35
$ws.Range("AI151").Value = 12
$ws.Range("AJ151").Value = 1.424
$ws.Range("AK151").Value = 0.318
$ws.Range("AL151").Value = 4.6
$ws.Range("AM151").Value = 2

# Row 157
$ws.Range("E157").Value = 6
$ws.Range("F157").Value = 6
$ws.Range("J157").Value = 12
$ws.Range("K157").Value = 9
$ws.Range("L157").Value = 2
$ws.Range("M157").Value = 2
$ws.Range("N157").Value = 3
$ws.Range("O157").Value = 13
$ws.Range("P157").Value = 1
$ws.Range("R157").Value = 1.5
$ws.Range("S157").Value = 43
$ws.Range("AB157").Value = 47
$ws.Range("AC157").Value = 186
$ws.Range("AG157").Value = 35
$ws.Range("AH157").Value = 26
$ws.Range("AI157").Value = 13
$ws.Range("AJ157").Value = 1
$ws.Range("AK157").Value = 0.267
$ws.Range("AL157").Value = 9.800000000000001
$ws.Range("AM157").Value = 4.33

# Row 170
$ws.Range("E170").Value = 20
$ws.Range("F170").Value = 20
$ws.Range("H170").Value = 10
$ws.Range("J170").Value = 103.2
$ws.Range("K170").Value = 91
$ws.Range("L170").Value = 56
$ws.Range("M170").Value = 55
$ws.Range("N170").Value = 43
$ws.Range("O170").Value = 130
$ws.Range("R170").Value = 4.77
$ws.Range("S170").Value = 391
$ws.Range("AB170").Value = 438
$ws.Range("AC170").Value = 1824
$ws.Range("AJ170").Value = 1.293
$ws.Range("AK170").Value = 0.296
$ws.Range("AL170").Value = 11.3
$ws.Range("AM170").Value = 3.02

# Row 190
$ws.Range("E190").Value = 19
$ws.Range("J190").Value = 39.2
$ws.Range("K190").Value = 30
$ws.Range("L190").Value = 12
$ws.Range("M190").Value = 11
$ws.Range("N190").Value = 12
$ws.Range("P190").Value = 5
$ws.Range("R190").Value = 2.5
$ws.Range("S190").Value = 144
$ws.Range("T190").Value = 2
$ws.Range("AB190").Value = 157
$ws.Range("AC190").Value = 601
$ws.Range("AD190").Value = 64
$ws.Range("AF190").Value = 7
$ws.Range("AH190").Value = 27
$ws.Range("AJ190").Value = 1.059
$ws.Range("AK190").Value = 0.207
$ws.Range("AL190").Value = 4.3
$ws.Range("AM190").Value = 1.58

# Row 193
$ws.Range("E193").Value = 43
$ws.Range("J193").Value = 40.1
$ws.Range("N193").Value = 13
$ws.Range("O193").Value = 46
$ws.Range("R193").Value = 3.12
$ws.Range("S193").Value = 147
$ws.Range("AB193").Value = 166
$ws.Range("AC193").Value = 658
$ws.Range("AG193").Value = 24
$ws.Range("AH193").Value = 25
$ws.Range("AJ193").Value = 1.14
$ws.Range("AK193").Value = 0.287
$ws.Range("AM193").Value = 3.54

# Row 195
$ws.Range("E195").Value = 17
$ws.Range("F195").Value = 17
$ws.Range("J195").Value = 85
$ws.Range("K195").Value = 93
$ws.Range("L195").Value = 58
$ws.Range("M195").Value = 54
$ws.Range("N195").Value = 36
$ws.Range("O195").Value = 67
$ws.Range("R195").Value = 5.72
$ws.Range("S195").Value = 334
$ws.Range("T195").Value = 13
$ws.Range("W195").Value = 9
$ws.Range("Y195").Value = 15
$ws.Range("AB195").Value = 381
$ws.Range("AC195").Value = 1478
$ws.Range("AD195").Value = 63
$ws.Range("AG195").Value = 46
$ws.Range("AJ195").Value = 1.518
$ws.Range("AK195").Value = 0.302
$ws.Range("AL195").Value = 7.1
$ws.Range("AM195").Value = 1.86

# Row 199
$ws.Range("E199").Value = 19
$ws.Range("F199").Value = 19
$ws.Range("J199").Value = 103
$ws.Range("K199").Value = 113
$ws.Range("L199").Value = 54
$ws.Range("M199").Value = 52
$ws.Range("N199").Value = 26
$ws.Range("O199").Value = 63
$ws.Range("P199").Value = 20
$ws.Range("Q199").Value = 5
$ws.Range("R199").Value = 4.54
$ws.Range("S199").Value = 412
$ws.Range("T199").Value = 15
$ws.Range("AB199").Value = 445
$ws.Range("AC199").Value = 1665
$ws.Range("AJ199").Value = 1.35
$ws.Range("AK199").Value = 0.281
$ws.Range("AL199").Value = 5.5
$ws.Range("AM199").Value = 2.42

# Row 207
$ws.Range("E207").Value = 21
$ws.Range("F207").Value = 21
$ws.Range("H207").Value = 7
$ws.Range("J207").Value = 124.2
$ws.Range("K207").Value = 85
$ws.Range("L207").Value = 49
$ws.Range("M207").Value = 43
$ws.Range("N207").Value = 46
$ws.Range("O207").Value = 139
$ws.Range("R207").Value = 3.1
$ws.Range("S207").Value = 452
$ws.Range("T207").Value = 20
$ws.Range("AB207").Value = 505
$ws.Range("AC207").Value = 2041
$ws.Range("AJ207").Value = 1.051
$ws.Range("AK207").Value = 0.234
$ws.Range("AL207").Value = 10
$ws.Range("AM207").Value = 3.02

# Row 218
$ws.Range("E218").Value = 7
$ws.Range("F218").Value = 7
$ws.Range("H218").Value = 4
$ws.Range("J218").Value = 34.1
$ws.Range("K218").Value = 50
$ws.Range("L218").Value = 25
$ws.Range("M218").Value = 23
$ws.Range("O218").Value = 20
$ws.Range("P218").Value = 7
$ws.Range("R218").Value = 6.03
$ws.Range("S218").Value = 149
$ws.Range("T218").Value = 12
$ws.Range("X218").Value = 1
$ws.Range("Y218").Value = 3
$ws.Range("AB218").Value = 158
$ws.Range("AC218").Value = 580
$ws.Range("AF218").Value = 8
$ws.Range("AH218").Value = 19
$ws.Range("AI218").Value = 6
$ws.Range("AJ218").Value = 1.66
$ws.Range("AK218").Value = 0.35
$ws.Range("AL218").Value = 5.2
$ws.Range("AM218").Value = 2.86

# Row 226
$ws.Range("E226").Value = 20
$ws.Range("F226").Value = 20
$ws.Range("G226").Value = 7
$ws.Range("J226").Value = 116.2
$ws.Range("K226").Value = 97
$ws.Range("L226").Value = 54
$ws.Range("M226").Value = 52
$ws.Range("N226").Value = 36
$ws.Range("O226").Value = 112
$ws.Range("P226").Value = 15
$ws.Range("R226").Value = 4.01
$ws.Range("S226").Value = 433
$ws.Range("T226").Value = 21
$ws.Range("V226").Value = 2
$ws.Range("Y226").Value = 11
$ws.Range("AB226").Value = 478
$ws.Range("AC226").Value = 1892
$ws.Range("AE226").Value = 16
$ws.Range("AH226").Value = 26
$ws.Range("AJ226").Value = 1.14
$ws.Range("AK226").Value = 0.265
$ws.Range("AL226").Value = 8.6
$ws.Range("AM226").Value = 3.11

# Row 234
$ws.Range("E234").Value = 20
$ws.Range("F234").Value = 20
$ws.Range("J234").Value = 128
$ws.Range("K234").Value = 86
$ws.Range("L234").Value = 36
$ws.Range("M234").Value = 34
$ws.Range("O234").Value = 164
$ws.Range("P234").Value = 14
$ws.Range("Q234").Value = 5
$ws.Range("R234").Value = 2.39
$ws.Range("S234").Value = 459
$ws.Range("T234").Value = 19
$ws.Range("W234").Value = 8
$ws.Range("Y234").Value = 2
$ws.Range("Z234").Value = 3
$ws.Range("AB234").Value = 491
$ws.Range("AC234").Value = 2013
$ws.Range("AE234").Value = 14
$ws.Range("AJ234").Value = 0.875
$ws.Range("AK234").Value = 0.255
$ws.Range("AL234").Value = 11.5
$ws.Range("AM234").Value = 6.31

# Row 240
$ws.Range("E240").Value = 20
$ws.Range("F240").Value = 20
$ws.Range("J240").Value = 107.2
$ws.Range("K240").Value = 108
$ws.Range("L240").Value = 57
$ws.Range("M240").Value = 53
$ws.Range("N240").Value = 44
$ws.Range("O240").Value = 86
$ws.Range("P240").Value = 18
$ws.Range("R240").Value = 4.43
$ws.Range("S240").Value = 424
$ws.Range("AB240").Value = 475
$ws.Range("AC240").Value = 1841
$ws.Range("AF240").Value = 11
$ws.Range("AH240").Value = 25
$ws.Range("AI240").Value = 11
$ws.Range("AJ240").Value = 1.412
$ws.Range("AK240").Value = 0.278
$ws.Range("AL240").Value = 7.2
$ws.Range("AM240").Value = 1.95

